$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
